# fix random function initialisation: ignore None and empty strings
#
# Adds a new "static_one" test-fixture row (row 14) to Sheet1, mirroring the
# existing rows that exercise the numpy.random "choice" distribution, and
# restores the expected sheet/selection state (Sheet1 active, Sheet2 not).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Row 14 on Sheet1: copy number formats from the analogous row 2 cells so
# the new cells line up with the existing style palette (no new styles) ---
$ws1.Range("H2").Copy()
$ws1.Range("H14").PasteSpecial(-4122)

$ws1.Range("I2").Copy()
$ws1.Range("I14").PasteSpecial(-4122)

$ws1.Range("J2").Copy()
$ws1.Range("J14").PasteSpecial(-4122)

$ws1.Range("K2").Copy()
$ws1.Range("K14").PasteSpecial(-4122)

$ws1.Range("L2").Copy()
$ws1.Range("L14").PasteSpecial(-4122)

# --- Fill in the row 14 values ---
$ws1.Range("A14").Value = "static_one"
$ws1.Range("C14").Value = "numpy.random"
$ws1.Range("D14").Value = "choice"
$ws1.Range("E14").Value = 1
$ws1.Range("H14").Value = "litres"
$ws1.Range("I14").Value = 39814
$ws1.Range("J14").Value = 39904
$ws1.Range("K14").Value = 0.1
$ws1.Range("L14").Value = 39814
$ws1.Range("M14").Value = "test var 1"

# --- Make Sheet1 the active/selected sheet with A14:M14 highlighted, and
# leave Sheet2's selection untouched (F3) but no longer the active tab ---
$null = $ws1.Activate()
$null = $ws1.Range("A14:M14").Select()
